# Update "main pages" metadata on the ConsentComponentType ValueSet workbook:
#  - Translate the Publisher name from German to English
#  - Translate the Contact name from German to English
#  - Fill in the previously-empty Description value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Publisher (row 9): German legal name -> English translation
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Contact (row 10): same translation, keeping the trailing URL
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Description (row 12): previously empty, now populated
$ws.Range("B12").Value = "Types of consent components relevant to gICS for differentiation as a search criterion"
